# Paper revision June 2021
#
# The crops dataset dropped several minor/"not elsewhere classified" crop
# categories. Remove the corresponding rows (matched by the FAO name in
# column A) from the sheet; everything else stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$namesToRemove = @(
    "BUCKWHEAT, MILLET, CANARY SEED",
    "CITRUS FRUIT NES",
    "FRUIT NES FRESH",
    "GRAPES",
    "LINSEED",
    "LUPINS",
    "OILSEEDS NES",
    "OLIVES",
    "ROOTS, TUBERS NES",
    "SORGHUM"
)

$lastRow = $ws.UsedRange.Rows.Count

# Walk bottom-to-top so deleting a row never invalidates the row numbers
# still left to visit.
for ($r = $lastRow; $r -ge 1; $r--) {
    $name = $ws.Cells.Item($r, 1).Value()
    if ($namesToRemove -contains $name) {
        $ws.Rows.Item($r).Delete()
    }
}
